# Re-sort the worksheet tabs into the reverse of their current order:
#   before: 2021-Q4, 2022-Q1, 2022-Q2, 总计
#   after : 总计, 2022-Q2, 2022-Q1, 2021-Q4
$wb = $excel.ActiveWorkbook

$count = $wb.Worksheets.Count
for ($i = 2; $i -le $count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Move($wb.Worksheets.Item(1))
}
